$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 378519.7
$ws.Range("J17").Value = 412028.16
$ws.Range("L17").Value = 1236084.48
$ws.Range("N17").Value = -1236420.48
$ws.Range("H44").Value = 3050
$ws.Range("J44").Value = 3050
$ws.Range("L44").Value = 3050
$ws.Range("N44").Value = -3974
$ws.Range("H74").Value = 74105.3
$ws.Range("I74").Value = 108679.1
$ws.Range("J74").Value = 24714.143
$ws.Range("K74").Value = 108679.1
$ws.Range("L74").Value = 24714.143
$ws.Range("M74").Value = -107743.1
$ws.Range("N74").Value = -26586.143
$ws.Range("H77").Value = 74105.3
$ws.Range("I77").Value = 108679.1
$ws.Range("J77").Value = 24714.143
$ws.Range("K77").Value = 543395.5
$ws.Range("L77").Value = 123570.715
$ws.Range("M77").Value = -538715.5
$ws.Range("N77").Value = -132930.715
$ws.Range("H136").Value = 81619.8
$ws.Range("J136").Value = 81619.8
$ws.Range("L136").Value = 81619.8
$ws.Range("N136").Value = -91819.8
$ws.Range("H137").Value = 1400.5151
$ws.Range("I137").Value = 1376.6
$ws.Range("K137").Value = 4129.799999999999
$ws.Range("M137").Value = -1579.799999999999
$ws.Range("H138").Value = 4699.7095
$ws.Range("J138").Value = 6634.263
$ws.Range("L138").Value = 19902.789
$ws.Range("N138").Value = -30182.789

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6637.346
$ws.Range("I45").Value = 14884.875
$ws.Range("J45").Value = 2971.7778
$ws.Range("K45").Value = 14884.875
$ws.Range("L45").Value = 2971.7778
$ws.Range("M45").Value = -14507.875
$ws.Range("N45").Value = -3725.7778
$ws.Range("H74").Value = 5068.6665
$ws.Range("I74").Value = 978.35297
$ws.Range("K74").Value = 978.35297
$ws.Range("M74").Value = -104.35297
$ws.Range("H77").Value = 5068.6665
$ws.Range("I77").Value = 978.35297
$ws.Range("K77").Value = 4891.76485
$ws.Range("M77").Value = -523.7648500000005
$ws.Range("H102").Value = 5353.2
$ws.Range("I102").Value = 5441.5
$ws.Range("K102").Value = 5441.5
$ws.Range("M102").Value = -3819.5
$ws.Range("H122").Value = 2147.7693
$ws.Range("I122").Value = 1501.2222
$ws.Range("K122").Value = 4503.6666
$ws.Range("M122").Value = -2053.6666
$ws.Range("H132").Value = 5944
$ws.Range("I132").Value = 5424.625
$ws.Range("J132").Value = 6982.75
$ws.Range("K132").Value = 16273.875
$ws.Range("L132").Value = 20948.25
$ws.Range("M132").Value = -13743.875
$ws.Range("N132").Value = -26008.25
$ws.Range("H135").Value = 59102.5
$ws.Range("J135").Value = 59102.5
$ws.Range("L135").Value = 59102.5
$ws.Range("N135").Value = -69242.5
$ws.Range("H138").Value = 76754
$ws.Range("J138").Value = 76754
$ws.Range("L138").Value = 76754
$ws.Range("N138").Value = -87034

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2851.4517
$ws.Range("I134").Value = 2489.5417
$ws.Range("J134").Value = 4092.2856
$ws.Range("K134").Value = 7468.625100000001
$ws.Range("L134").Value = 12276.8568
$ws.Range("M134").Value = -4933.625100000001
$ws.Range("N134").Value = -17346.8568
$ws.Range("H139").Value = 70330
$ws.Range("J139").Value = 68396
$ws.Range("L139").Value = 68396
$ws.Range("N139").Value = -78676

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 10000
$ws.Range("I45").Value = 10000
$ws.Range("K45").Value = 10000
$ws.Range("M45").Value = -9407

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 26799.334
$ws.Range("I82").Value = 18499.166
$ws.Range("K82").Value = 55497.49800000001
$ws.Range("M82").Value = -55091.49800000001
$ws.Range("H85").Value = 26799.334
$ws.Range("I85").Value = 18499.166
$ws.Range("K85").Value = 55497.49800000001
$ws.Range("M85").Value = -54093.49800000001
$ws.Range("H121").Value = 6210.8184
$ws.Range("I121").Value = 4505.8887
$ws.Range("J121").Value = 7391.154
$ws.Range("K121").Value = 13517.6661
$ws.Range("L121").Value = 22173.462
$ws.Range("M121").Value = -12207.6661
$ws.Range("N121").Value = -24793.462
$ws.Range("H131").Value = 38924.668
$ws.Range("I131").Value = 500490
$ws.Range("J131").Value = 1999.44
$ws.Range("K131").Value = 1501470
$ws.Range("L131").Value = 5998.32
$ws.Range("M131").Value = -1496430
$ws.Range("N131").Value = -16078.32

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 30750
$ws.Range("J47").Value = 30750
$ws.Range("L47").Value = 30750
$ws.Range("N47").Value = -31886
$ws.Range("H48").Value = 22500
$ws.Range("I48").Value = 15000
$ws.Range("J48").Value = 30000
$ws.Range("K48").Value = 15000
$ws.Range("L48").Value = 30000
$ws.Range("M48").Value = -14515
$ws.Range("N48").Value = -30970
$ws.Range("H113").Value = 1838.7778
$ws.Range("I113").Value = 1838.7778
$ws.Range("K113").Value = 1838.7778
$ws.Range("M113").Value = 331.2221999999999
$ws.Range("H122").Value = 3227.25
$ws.Range("I122").Value = 3000.5
$ws.Range("J122").Value = 3454
$ws.Range("K122").Value = 9001.5
$ws.Range("L122").Value = 10362
$ws.Range("M122").Value = -6551.5
$ws.Range("N122").Value = -15262
$ws.Range("H126").Value = 14180.607
$ws.Range("I126").Value = 22353
$ws.Range("K126").Value = 67059
$ws.Range("M126").Value = -64589
$ws.Range("H132").Value = 558519.4
$ws.Range("I132").Value = 717167.8
$ws.Range("K132").Value = 2151503.4
$ws.Range("M132").Value = -2148973.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 26140
$ws.Range("J4").Value = 14210
$ws.Range("L4").Value = 14210
$ws.Range("N4").Value = -14436
$ws.Range("H22").Value = 1042.8572
$ws.Range("I22").Value = 456.25
$ws.Range("J22").Value = 1277.5
$ws.Range("K22").Value = 456.25
$ws.Range("L22").Value = 1277.5
$ws.Range("M22").Value = -161.25
$ws.Range("N22").Value = -1867.5
$ws.Range("H27").Value = 1042.8572
$ws.Range("I27").Value = 456.25
$ws.Range("J27").Value = 1277.5
$ws.Range("K27").Value = 456.25
$ws.Range("L27").Value = 1277.5
$ws.Range("M27").Value = -349.25
$ws.Range("N27").Value = -1491.5
$ws.Range("H28").Value = 26140
$ws.Range("J28").Value = 14210
$ws.Range("L28").Value = 14210
$ws.Range("N28").Value = -14674
$ws.Range("H37").Value = 26140
$ws.Range("J37").Value = 14210
$ws.Range("L37").Value = 14210
$ws.Range("N37").Value = -14424
$ws.Range("H136").Value = 4544.1665
$ws.Range("I136").Value = 3703
$ws.Range("J136").Value = 8750
$ws.Range("K136").Value = 11109
$ws.Range("L136").Value = 26250
$ws.Range("M136").Value = -8559
$ws.Range("N136").Value = -31350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 262493.75
$ws.Range("I52").Value = 343325
$ws.Range("K52").Value = 343325
$ws.Range("M52").Value = -343099
$ws.Range("H113").Value = 723.1053000000001
$ws.Range("I113").Value = 683.44446
$ws.Range("J113").Value = 758.8
$ws.Range("K113").Value = 2050.33338
$ws.Range("L113").Value = 2276.4
$ws.Range("M113").Value = 119.66662
$ws.Range("N113").Value = -6616.4
$ws.Range("H137").Value = 125924
$ws.Range("J137").Value = 125924
$ws.Range("L137").Value = 125924
$ws.Range("N137").Value = -136124
$ws.Range("H141").Value = 165175
$ws.Range("J141").Value = 165175
$ws.Range("L141").Value = 165175
$ws.Range("N141").Value = -175535
